$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.395.60"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "'1.878.59"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'0.7169"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "'243.58"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.07928"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "'0.3144"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").Value = "'24.99"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").Value = "'0.08138"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").Value = "'1.884.01"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "'95.59"
$ws.Range("E13").Value = "  +4.74%  "
$ws.Range("D14").Value = "'5.238"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "'6.401"
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("D17").Value = "'0.000008397"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "'29.404.79"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'252.01"
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("D20").Value = "'13.36"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "'2.141.44"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'7.656"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "'9.071"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'162.21"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "'4.414"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'4.296"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").Value = "'1.217"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'0.05329"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'1.946"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "'0.7583"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("D36").Value = "'1.177"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'0.01894"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "'1.269.26"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").Value = "'2.762"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").Value = "'6.391"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "'112.04"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'0.9057"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").Value = "'74.27"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'2.035.14"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "'1.810"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "'9.510"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").Value = "'0.4349"
$ws.Range("E51").Value = "  +0.33%  "
